$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New BOM row 20: P1 (Würth touch-sensor pad component) ---
# Values are entered in this order so the shared-string table grows in the
# same sequence as the authored workbook (C -> F -> E -> D).
$ws.Range("C20").Value = "P1"
$ws.Range("F20").Value = "we - 3029040030025"
$ws.Range("E20").Value = "Würth - 3029040030025"
$ws.Range("D20").Value = "WE-SMGS_3029040030025"
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 1

# Formatting: reuse the formatting of row 12 (same banded style) for the
# plain cells, and build the text-formatted / hyperlink-formatted cells for
# the footprint/description/comment columns.
$ws.Range("A12").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("B12").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("G12").Copy()
$ws.Range("G20").PasteSpecial(-4122)

$ws.Range("D10").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").NumberFormat = "@"

$ws.Range("E10").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").NumberFormat = "@"

$ws.Range("F12").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.we-online.com/catalog/en/element/3029040030025") | Out-Null

$excel.CutCopyMode = 0

# Selection moves on, matching where the author left off after the edit.
$ws.Range("D24").Select()
